$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove the two extra test-data rows (rows 3 & 4) from advanced_search,
# leaving only the header row and a single data row.
$ws1.Rows.Item(4).Delete() | Out-Null
$ws1.Rows.Item(3).Delete() | Out-Null

# Update the remaining data row with the new drug/param values.
$ws1.Range("C2").Value = "d=C2039&loc=0&rl=2"
$ws1.Range("B2").Value = "Bevacizumab"

# advanced_search becomes the active sheet / selected tab, with A7 selected.
# (advanced_search_negative keeps its own D6 selection untouched, and loses
# tabSelected automatically since it is no longer the active sheet.)
$ws1.Activate()
$ws1.Range("A7").Select() | Out-Null

Write-Output "done"
